$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Explanation" column header and value
$ws.Range("H1").Value = "Explanation"
$ws.Range("H2").Value = "Explanation to why it's right"

# Set column width for new column H
# (Target OOXML width is 28.88671875; due to the runtime's internal
# pixel-based quantization of column widths, a ColumnWidth of 28 lands
# in the closest achievable bucket.)
$ws.Columns.Item(8).ColumnWidth = 28

# Update selection to match target (I9)
$ws.Range("I9").Select()
